$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The difference_percentage column holds plain text like "1.97%". Excel's
# normal Value assignment would auto-coerce such strings into numeric
# percentages, so mark the cells as Text first, then restore a plain
# (unstyled) look afterwards.
$ws.Range("E2:E4").NumberFormat = "@"

# Update row 2
$ws.Range("A2").Value = "No significant differences detected between register_clicked_register.png_20240807-092154.png and register_clicked_register.png_20240807-095206.png."
$ws.Range("B2").Value = "Success"
$ws.Range("C2").Value = "2024-08-07 09:52:25"
$ws.Range("D2").Value = "Master"
$ws.Range("E2").Value = "1.97%"
$ws.Range("F2").Value = "register_clicked_register.png_20240807-092154.png"
$ws.Range("G2").Value = "register_clicked_register.png_20240807-095206.png"

# Update row 3
$ws.Range("A3").Value = "No significant differences detected between register_filled_form.png_20240807-092157.png and register_filled_form.png_20240807-095209.png."
$ws.Range("B3").Value = "Success"
$ws.Range("C3").Value = "2024-08-07 09:52:25"
$ws.Range("D3").Value = "Master"
$ws.Range("E3").Value = "0.33%"
$ws.Range("F3").Value = "register_filled_form.png_20240807-092157.png"
$ws.Range("G3").Value = "register_filled_form.png_20240807-095209.png"

# Add new row 4
$ws.Range("A4").Value = "No significant differences detected between register_submitted.png_20240807-092211.png and register_submitted.png_20240807-095222.png."
$ws.Range("B4").Value = "Success"
$ws.Range("C4").Value = "2024-08-07 09:52:25"
$ws.Range("D4").Value = "Master"
$ws.Range("E4").Value = "0.09%"
$ws.Range("F4").Value = "register_submitted.png_20240807-092211.png"
$ws.Range("G4").Value = "register_submitted.png_20240807-095222.png"

# Restore the default (unstyled) look on the percentage cells now that the
# text values are locked in, matching the rest of the data rows.
$ws.Range("E2:E4").Style = "Normal"
